$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (prices/volumes/order) as scraped on Fri Jan 27 15:59:55 UTC 2023.
# Each target cell is forced to Text (NumberFormat "@") before the write so numeric-looking
# strings like "304.39" or "-0.52%" are stored as text (matching the sheet's inlineStr cells)
# instead of being auto-converted to numbers; ClearFormats() afterwards restores the original
# (default) cell style so only the value itself changes.

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "304.39"
Set-TextValue "E2" "-0.52%"
Set-TextValue "D3" "35.83"
Set-TextValue "E3" "0.41%"
Set-TextValue "D4" "5.057"
Set-TextValue "E4" "-0.42%"
Set-TextValue "D5" "0.08005"
Set-TextValue "E5" "-0.70%"
Set-TextValue "D6" "1.862"
Set-TextValue "E6" "-4.00%"
Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "7.772"
Set-TextValue "E7" "-0.77%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9259"
Set-TextValue "E8" "-1.38%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1290"
Set-TextValue "E9" "-6.53%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1894"
Set-TextValue "E10" "0.03%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.09026"
Set-TextValue "E11" "-1.53%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03446"
Set-TextValue "E12" "-1.83%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09878"
Set-TextValue "E13" "-0.07%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001405"
Set-TextValue "E14" "-3.61%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.006238"
Set-TextValue "E15" "-6.92%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.834"
Set-TextValue "E16" "5.80%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.123"
Set-TextValue "E17" "-1.31%"
Set-TextValue "D19" "0.3408"
Set-TextValue "E19" "-0.58%"
Set-TextValue "D20" "0.1335"
Set-TextValue "E20" "-0.90%"
Set-TextValue "D21" "4.827"
Set-TextValue "E21" "-8.12%"
Set-TextValue "D22" "0.2413"
Set-TextValue "E22" "-5.08%"
Set-TextValue "D23" "0.04362"
Set-TextValue "E23" "-1.19%"
Set-TextValue "D24" "0.001231"
Set-TextValue "E24" "-0.65%"
Set-TextValue "D25" "0.004820"
Set-TextValue "E25" "2.02%"
Set-TextValue "D27" "0.0001302"
Set-TextValue "E27" "-0.26%"
Set-TextValue "E28" "41.57%"
Set-TextValue "D39" "0.01960"
Set-TextValue "E39" "-1.93%"
Set-TextValue "D40" "0.05115"
Set-TextValue "E40" "0.32%"
Set-TextValue "D41" "0.007503"
Set-TextValue "E41" "-1.82%"
Set-TextValue "D42" "0.01012"
Set-TextValue "E42" "-9.58%"
Set-TextValue "D43" "0.1351"
Set-TextValue "E43" "-1.46%"
Set-TextValue "D44" "0.002113"
Set-TextValue "E44" "0.22%"
Set-TextValue "D45" "0.009865"
Set-TextValue "E45" "-12.85%"
Set-TextValue "D46" "0.00006184"
Set-TextValue "E46" "-2.86%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "-0.28%"
Set-TextValue "D48" "64.95"
Set-TextValue "E48" "-0.42%"
Set-TextValue "D49" "0.001252"
Set-TextValue "E49" "-22.11%"
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "-0.28%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "-0.28%"
